$d = $word.ActiveDocument

# Locate the target paragraph (last one, "31/05/2016 : ...").
$p = $d.Paragraphs.Item(9)

# Position right after "...plus proches voisins et " -- this is exactly where
# the existing "_GoBack" bookmark currently sits (between that run and the
# run containing "des distances du point à ses ").
$f = $d.Content
$f.Find.Execute("plus proches voisins et ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkOldPos = $f.End

# Relocate the "_GoBack" bookmark out of the way (to the very start of the
# paragraph) before we touch the text around its old position, so it does
# not interfere with -- or get clobbered by -- the edits below.
$moveAway = $d.Range($p.Range.Start, $p.Range.Start)
$d.Bookmarks.Add("_GoBack", $moveAway)

# Append the new trailing sentence at the very end of the paragraph (right
# after "ppv", before the paragraph mark).
$endOfPara = $p.Range.End - 1
$tail = $d.Range($endOfPara, $endOfPara)
$tail.InsertAfter(". Modification de la gdb ajout nom et type par défaut")

# Re-read the paragraph end now that new text has been appended.
$endOfPara = $p.Range.End - 1

# Work around an edge-case bug: placing a bookmark exactly at the
# paragraph-end position (collapsed range immediately before the paragraph
# mark) makes it jump to the start of the document. Insert a one-character
# placeholder there, anchor the bookmark just before it, then remove the
# placeholder again.
$placeholderPos = $endOfPara
$ph = $d.Range($placeholderPos, $placeholderPos)
$ph.InsertAfter("X")

$bmRange = $d.Range($placeholderPos, $placeholderPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$phRange = $d.Range($placeholderPos, $placeholderPos + 1)
$phRange.Text = ""
